$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 1
    4  = 1
    5  = 2
    6  = 0
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 5
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    23 = 3
    24 = 1
    25 = 3
    26 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
